$wb = $excel.ActiveWorkbook

# --- Sheet1 "TestCases": update selection (tabSelected moves to sheet2) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A2").Select()

# --- Sheet2 "Test_A": trim to A1:C3, replace Age/Subject/Runmode/Result/Error
#     columns (D:G) with a single "Div" column, rewrite the data rows, rename
#     the sheet, and make it the active tab with C3 selected. ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("D1:G4").EntireColumn.Delete()
$ws2.Range("A4:C4").EntireRow.Delete()

$ws2.Range("A2").Value = "admin"
$ws2.Range("C1").Value = "Div"
$ws2.Range("C2").Value = "!st"
$ws2.Range("A3").Value = "fadmin"
$ws2.Range("B2").Value = "csm10002"
$ws2.Range("B3").Value = "csm10003"
$ws2.Range("C3").Value = "2nd"

$ws2.Name = "verifyISPUserSuccessfulLogin"
$ws2.Activate()
$ws2.Range("C3").Select()
